# Updated symbol list on Fri Feb 10 22:35:49 UTC 2023 with GitHub Actions
# Refreshes Coin / Link / Price / Volume(1h) cells on the cryptos sheet,
# mirroring the latest coinranking.com snapshot pulled by the scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = "D2"; Value = "306.61" }
    @{ Cell = "E2"; Value = "-0.42%" }
    @{ Cell = "D3"; Value = "40.31" }
    @{ Cell = "E3"; Value = "1.12%" }
    @{ Cell = "D4"; Value = "5.104" }
    @{ Cell = "E4"; Value = "0.70%" }
    @{ Cell = "D5"; Value = "0.07586" }
    @{ Cell = "E5"; Value = "-2.59%" }
    @{ Cell = "B6"; Value = "GateToken" }
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
    @{ Cell = "D6"; Value = "4.277" }
    @{ Cell = "E6"; Value = "-0.99%" }
    @{ Cell = "B7"; Value = "FTXToken" }
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" }
    @{ Cell = "D7"; Value = "1.608" }
    @{ Cell = "E7"; Value = "-2.82%" }
    @{ Cell = "B8"; Value = "BTSEToken" }
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Cell = "D8"; Value = "2.447" }
    @{ Cell = "E8"; Value = "-4.41%" }
    @{ Cell = "B9"; Value = "MXToken" }
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Cell = "D9"; Value = "0.9053" }
    @{ Cell = "E9"; Value = "-1.47%" }
    @{ Cell = "B10"; Value = "LiechtensteinCryptoassetsExchange" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" }
    @{ Cell = "D10"; Value = "0.1011" }
    @{ Cell = "E10"; Value = "1.98%" }
    @{ Cell = "B11"; Value = "WazirX" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Cell = "D11"; Value = "0.1754" }
    @{ Cell = "E11"; Value = "0.80%" }
    @{ Cell = "B12"; Value = "MandalaExchangeToken" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "D12"; Value = "0.09142" }
    @{ Cell = "E12"; Value = "1.67%" }
    @{ Cell = "B13"; Value = "BitrueCoin" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D13"; Value = "0.04186" }
    @{ Cell = "E13"; Value = "-4.73%" }
    @{ Cell = "B14"; Value = "BitMartToken" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "D14"; Value = "0.1055" }
    @{ Cell = "E14"; Value = "-0.52%" }
    @{ Cell = "B15"; Value = "BitForexToken" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "D15"; Value = "0.001233" }
    @{ Cell = "E15"; Value = "-1.81%" }
    @{ Cell = "B16"; Value = "TigerCash" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D16"; Value = "0.005836" }
    @{ Cell = "E16"; Value = "3.22%" }
    @{ Cell = "B17"; Value = "LEO" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D17"; Value = "3.349" }
    @{ Cell = "E17"; Value = "-0.51%" }
    @{ Cell = "D18"; Value = "0.3273" }
    @{ Cell = "D19"; Value = "6.661" }
    @{ Cell = "E19"; Value = "-5.72%" }
    @{ Cell = "D20"; Value = "0.1357" }
    @{ Cell = "E20"; Value = "-0.59%" }
    @{ Cell = "D21"; Value = "0.2732" }
    @{ Cell = "E21"; Value = "2.73%" }
    @{ Cell = "D22"; Value = "0.04184" }
    @{ Cell = "E22"; Value = "0.94%" }
    @{ Cell = "E23"; Value = "1.81%" }
    @{ Cell = "D24"; Value = "0.004055" }
    @{ Cell = "E24"; Value = "-0.91%" }
    @{ Cell = "D25"; Value = "0.0001304" }
    @{ Cell = "E25"; Value = "6.52%" }
    @{ Cell = "D26"; Value = "0.0003016" }
    @{ Cell = "E26"; Value = "0.76%" }
    @{ Cell = "D38"; Value = "0.02402" }
    @{ Cell = "E38"; Value = "0.56%" }
    @{ Cell = "D39"; Value = "0.05131" }
    @{ Cell = "E39"; Value = "-1.10%" }
    @{ Cell = "D40"; Value = "0.007788" }
    @{ Cell = "E40"; Value = "-1.99%" }
    @{ Cell = "D41"; Value = "0.1291" }
    @{ Cell = "E41"; Value = "-2.76%" }
    @{ Cell = "D42"; Value = "0.007055" }
    @{ Cell = "E42"; Value = "-0.70%" }
    @{ Cell = "D43"; Value = "0.001944" }
    @{ Cell = "E43"; Value = "-3.66%" }
    @{ Cell = "D44"; Value = "0.008465" }
    @{ Cell = "E44"; Value = "5.25%" }
    @{ Cell = "D45"; Value = "0.3323" }
    @{ Cell = "D46"; Value = "0.00006373" }
    @{ Cell = "E46"; Value = "-5.25%" }
    @{ Cell = "E47"; Value = "-0.17%" }
    @{ Cell = "B48"; Value = "BOLO" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo" }
    @{ Cell = "D48"; Value = "0.03166" }
    @{ Cell = "E48"; Value = "825.49%" }
    @{ Cell = "B49"; Value = "CoinbaseStockToken" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin" }
    @{ Cell = "D49"; Value = "0.004414" }
    @{ Cell = "E49"; Value = "7.12%" }
    @{ Cell = "D50"; Value = "0.00002107" }
    @{ Cell = "E50"; Value = "-0.17%" }
    @{ Cell = "D51"; Value = "0.0002007" }
    @{ Cell = "E51"; Value = "-0.17%" }
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    # Briefly force text storage so numeric-looking strings (prices,
    # percentage deltas) are written as text rather than being
    # coerced into Excel numbers, then drop the temporary format so
    # the cell keeps its original (default) styling.
    $cell.NumberFormat = "@"
    $cell.Value = $edit.Value
    $cell.ClearFormats()
}

